$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 3: the Intel driver row becomes the Realtek driver row (values
# that used to live on row 4), with the same number/text formatting.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 106
$ws.Range("D3").Value = 98.5

# ------------------------------------------------------------------
# Row 4: pick up the formatting from row 5 (the old "Totals:" row)
# before we touch it, then turn row 4 into the new Totals row.
# ------------------------------------------------------------------
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B5:C5").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "Totals:"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 106
$ws.Range("D4").Clear()

# ------------------------------------------------------------------
# Row 5: this was the old Totals row; it is now blank.
# ------------------------------------------------------------------
$ws.Range("A5:C5").Clear()

# ------------------------------------------------------------------
# Remove the blank row 10 (shifts the "Good Drivers" block - old
# rows 11-16 - up by one, so the header lands on row 10 and the
# column-header row lands on row 11).
# ------------------------------------------------------------------
$ws.Rows.Item(10).Delete()

# ------------------------------------------------------------------
# The three driver-detail rows (now at rows 12-14, previously 13-15)
# are no longer reported - clear their contents/formatting.
# ------------------------------------------------------------------
$ws.Range("A12:E14").Clear()
